# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Update the "Date" metadata value on the Metadata sheet.
# 2. Swap the "Mapping: RIM Mapping" (AK) and
#    "Mapping: Spécification métier vers l'extension ROR territorial
#    division" (AL) columns on the Elements sheet (header, data and
#    column widths).

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 : Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements : swap columns AK (37) and AL (38) ------------------------
$elements = $wb.Worksheets.Item("Elements")

$lastRow = $elements.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)

    $akValue = $akCell.Value2
    $alValue = $alCell.Value2

    # Only touch cells that actually need to change - rows where AK/AL
    # already hold an identical value (e.g. both blank) must stay
    # untouched so their underlying representation doesn't change.
    if ($akValue -ne $alValue) {
        $akCell.Value = $alValue
        $alCell.Value = $akValue
    }
}

# Swap the column widths too: AK used to be the narrow "bestFit" column
# (~24.98) and AL the wide one (~73.82); after swapping the content the
# widths follow the content to the other column.
$elements.Columns.Item(37).ColumnWidth = 73.81640625
$elements.Columns.Item(38).ColumnWidth = 24.98046875
